$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "name" column becomes "itemName"
$ws.Range("C1").Value = "itemName"

# Item display names translated to Korean
$ws.Range("C4").Value  = "하트 | HP 50+"
$ws.Range("C5").Value  = "이동속도 +"
$ws.Range("C6").Value  = "실드 |  방어력+"
$ws.Range("C7").Value  = "파워 +"
$ws.Range("C12").Value = "뿅망치"
$ws.Range("C13").Value = "낡은 검"
$ws.Range("C14").Value = "체인-쏘우"
$ws.Range("C15").Value = "다크 소드"
$ws.Range("C16").Value = "서리한"
$ws.Range("C17").Value = "철퇴"

# Mace price bump
$ws.Range("E17").Value = 600

# Remove stray K17 value
$ws.Range("K17").ClearContents()

# Update selection to reflect final cursor position
$ws.Range("K17").Select()
